$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 4189.4785
$ws.Range("I15").Value = 4189.4785
$ws.Range("K15").Value = 12568.4355
$ws.Range("M15").Value = -12399.4355

$ws.Range("H108").Value = 31282
$ws.Range("J108").Value = 31282
$ws.Range("L108").Value = 31282
$ws.Range("N108").Value = -38962

$ws.Range("H109").Value = 32940.668
$ws.Range("J109").Value = 32940.668
$ws.Range("L109").Value = 32940.668
$ws.Range("N109").Value = -35714.668

$ws.Range("H117").Value = 44371
$ws.Range("J117").Value = 44371
$ws.Range("L117").Value = 44371
$ws.Range("N117").Value = -53549

$ws.Range("H120").Value = 49714
$ws.Range("J120").Value = 49714
$ws.Range("L120").Value = 49714
$ws.Range("N120").Value = -59390

$ws.Range("H128").Value = 42581.75
$ws.Range("J128").Value = 42581.75
$ws.Range("L128").Value = 42581.75
$ws.Range("N128").Value = -52541.75

$ws.Range("H130").Value = 43340
$ws.Range("J130").Value = 43340
$ws.Range("L130").Value = 43340
$ws.Range("N130").Value = -53380

$ws.Range("H138").Value = 1921.2659
$ws.Range("I138").Value = 1147.035
$ws.Range("J138").Value = 3927.2273
$ws.Range("K138").Value = 3441.105
$ws.Range("L138").Value = 11781.6819
$ws.Range("M138").Value = 1698.895
$ws.Range("N138").Value = -22061.6819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10413.389
$ws.Range("I32").Value = 9029.365
$ws.Range("J32").Value = 20101.555
$ws.Range("K32").Value = 9029.365
$ws.Range("L32").Value = 20101.555
$ws.Range("M32").Value = -8742.365
$ws.Range("N32").Value = -20675.555

$ws.Range("H45").Value = 1792.4445
$ws.Range("I45").Value = 1556.6923
$ws.Range("K45").Value = 1556.6923
$ws.Range("M45").Value = -1179.6923

$ws.Range("H88").Value = 14818915
$ws.Range("I88").Value = 25002488
$ws.Range("J88").Value = 8552101
$ws.Range("K88").Value = 25002488
$ws.Range("L88").Value = 8552101
$ws.Range("M88").Value = -25002082
$ws.Range("N88").Value = -8552913

$ws.Range("H91").Value = 14818915
$ws.Range("I91").Value = 25002488
$ws.Range("J91").Value = 8552101
$ws.Range("K91").Value = 25002488
$ws.Range("L91").Value = 8552101
$ws.Range("M91").Value = -25001084
$ws.Range("N91").Value = -8554909

$ws.Range("H97").Value = 753.4828
$ws.Range("I97").Value = 611.73914
$ws.Range("J97").Value = 1296.8334
$ws.Range("K97").Value = 611.73914
$ws.Range("L97").Value = 1296.8334
$ws.Range("M97").Value = -115.73914
$ws.Range("N97").Value = -2288.8334

$ws.Range("H107").Value = 30794.6
$ws.Range("J107").Value = 30794.6
$ws.Range("L107").Value = 30794.6
$ws.Range("N107").Value = -38474.6

$ws.Range("H109").Value = 23427
$ws.Range("J109").Value = 23427
$ws.Range("L109").Value = 23427
$ws.Range("N109").Value = -26201

$ws.Range("H117").Value = 38562.8
$ws.Range("J117").Value = 38562.8
$ws.Range("L117").Value = 38562.8
$ws.Range("N117").Value = -47740.8

$ws.Range("H118").Value = 49356
$ws.Range("J118").Value = 49356
$ws.Range("L118").Value = 49356
$ws.Range("N118").Value = -52670

$ws.Range("H120").Value = 41240
$ws.Range("J120").Value = 41240
$ws.Range("L120").Value = 41240
$ws.Range("N120").Value = -50916

$ws.Range("H125").Value = 50698
$ws.Range("J125").Value = 50698
$ws.Range("L125").Value = 50698
$ws.Range("N125").Value = -60538

$ws.Range("H128").Value = 50421
$ws.Range("J128").Value = 50421
$ws.Range("L128").Value = 50421
$ws.Range("N128").Value = -60381

$ws.Range("H138").Value = 27500
$ws.Range("J138").Value = 27500
$ws.Range("L138").Value = 27500
$ws.Range("N138").Value = -37780

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2386.0908
$ws.Range("I86").Value = 2428.3572
$ws.Range("J86").Value = 2312.125
$ws.Range("K86").Value = 2428.3572
$ws.Range("L86").Value = 2312.125
$ws.Range("M86").Value = -1305.3572
$ws.Range("N86").Value = -4558.125

$ws.Range("H89").Value = 2386.0908
$ws.Range("I89").Value = 2428.3572
$ws.Range("J89").Value = 2312.125
$ws.Range("K89").Value = 12141.786
$ws.Range("L89").Value = 11560.625
$ws.Range("M89").Value = -6525.786
$ws.Range("N89").Value = -22792.625

$ws.Range("H119").Value = 42248.332
$ws.Range("J119").Value = 42248.332
$ws.Range("L119").Value = 42248.332
$ws.Range("N119").Value = -51924.332

$ws.Range("H120").Value = 46757
$ws.Range("J120").Value = 46757
$ws.Range("L120").Value = 46757
$ws.Range("N120").Value = -56433

$ws.Range("H125").Value = 50780
$ws.Range("J125").Value = 50780
$ws.Range("L125").Value = 50780
$ws.Range("N125").Value = -60620

$ws.Range("H126").Value = 43608
$ws.Range("J126").Value = 43608
$ws.Range("L126").Value = 43608
$ws.Range("N126").Value = -53488

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 46915.332
$ws.Range("J20").Value = 46915.332
$ws.Range("L20").Value = 46915.332
$ws.Range("N20").Value = -47387.332

$ws.Range("H30").Value = 46915.332
$ws.Range("J30").Value = 46915.332
$ws.Range("L30").Value = 46915.332
$ws.Range("N30").Value = -47097.332

$ws.Range("H62").Value = 2888.6
$ws.Range("I62").Value = 2665.7273
$ws.Range("J62").Value = 3501.5
$ws.Range("K62").Value = 2665.7273
$ws.Range("L62").Value = 3501.5
$ws.Range("M62").Value = -2041.7273
$ws.Range("N62").Value = -4749.5

$ws.Range("H65").Value = 2888.6
$ws.Range("I65").Value = 2665.7273
$ws.Range("J65").Value = 3501.5
$ws.Range("K65").Value = 13328.6365
$ws.Range("L65").Value = 17507.5
$ws.Range("M65").Value = -10208.6365
$ws.Range("N65").Value = -23747.5

$ws.Range("H86").Value = 5890.6
$ws.Range("I86").Value = 5915.2856
$ws.Range("J86").Value = 5833
$ws.Range("K86").Value = 5915.2856
$ws.Range("L86").Value = 5833
$ws.Range("M86").Value = -4792.2856
$ws.Range("N86").Value = -8079

$ws.Range("H89").Value = 5890.6
$ws.Range("I89").Value = 5915.2856
$ws.Range("J89").Value = 5833
$ws.Range("K89").Value = 29576.428
$ws.Range("L89").Value = 29165
$ws.Range("M89").Value = -23960.428
$ws.Range("N89").Value = -40397

$ws.Range("H116").Value = 42872
$ws.Range("J116").Value = 42872
$ws.Range("L116").Value = 42872
$ws.Range("N116").Value = -52050

$ws.Range("H128").Value = 46915.332
$ws.Range("J128").Value = 46915.332
$ws.Range("L128").Value = 46915.332
$ws.Range("N128").Value = -56875.332

$ws.Range("H132").Value = 17190
$ws.Range("I132").Value = 1111.1666
$ws.Range("J132").Value = 89544.75
$ws.Range("K132").Value = 3333.4998
$ws.Range("L132").Value = 268634.25
$ws.Range("M132").Value = -803.4998000000001
$ws.Range("N132").Value = -273694.25

$ws.Range("H134").Value = 275840.06
$ws.Range("I134").Value = 953.6842
$ws.Range("K134").Value = 2861.0526
$ws.Range("M134").Value = -326.0526

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5587.5
$ws.Range("I70").Value = 5400
$ws.Range("K70").Value = 5400
$ws.Range("M70").Value = -5130

$ws.Range("H73").Value = 5587.5
$ws.Range("I73").Value = 5400
$ws.Range("K73").Value = 5400
$ws.Range("M73").Value = -4464

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H121").Value = 43416
$ws.Range("J121").Value = 43416
$ws.Range("L121").Value = 43416
$ws.Range("N121").Value = -46910

$ws.Range("H127").Value = 50510.832
$ws.Range("J127").Value = 50510.832
$ws.Range("L127").Value = 50510.832
$ws.Range("N127").Value = -60430.832

$ws.Range("H132").Value = 2304.4602
$ws.Range("I132").Value = 1395.6957
$ws.Range("K132").Value = 4187.0871
$ws.Range("M132").Value = -1657.0871

$ws.Range("H136").Value = 1595.825
$ws.Range("I136").Value = 1027.2593
$ws.Range("J136").Value = 2776.6924
$ws.Range("K136").Value = 3081.7779
$ws.Range("L136").Value = 8330.0772
$ws.Range("M136").Value = -531.7779
$ws.Range("N136").Value = -13430.0772

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 45874
$ws.Range("J16").Value = 45874
$ws.Range("L16").Value = 45874
$ws.Range("N16").Value = -46458

$ws.Range("H119").Value = 200030740
$ws.Range("J119").Value = 200030740
$ws.Range("L119").Value = 200030740
$ws.Range("N119").Value = -200040416

$ws.Range("H120").Value = 40206
$ws.Range("J120").Value = 40206
$ws.Range("L120").Value = 40206
$ws.Range("N120").Value = -49882

$ws.Range("H124").Value = 28143
$ws.Range("J124").Value = 28143
$ws.Range("L124").Value = 28143
$ws.Range("N124").Value = -37963

$ws.Range("H128").Value = 44996
$ws.Range("J128").Value = 44996
$ws.Range("L128").Value = 44996
$ws.Range("N128").Value = -54956

$ws.Range("H132").Value = 4899.2
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 4899.2
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 14697.6
$ws.Range("N132").Value = -19757.6
$ws.Range("M132").ClearContents()

$ws.Range("H133").Value = 70091.39999999999
$ws.Range("J133").Value = 70091.39999999999
$ws.Range("L133").Value = 70091.39999999999
$ws.Range("N133").Value = -80211.39999999999
